# DeveloperGuide: update section of UndoRedoStack to VersionedAddressBook
#
# 1) The "datetimeFigureOut" date placeholder fields across the slide
#    master, every slide layout, and the notes master were re-cached by
#    PowerPoint (the deck was re-saved on a later date) — bump their
#    displayed text from "7/20/17" to "4/16/2018".
# 2) The UndoRedoStack shapes (the "UndoRedo / Stack" rectangle, its
#    arrow connector, and its "1" label textbox) on slide 1 are removed
#    now that the undo/redo mechanism has moved to VersionedAddressBook.

$p = $ppt.ActivePresentation

$oldDate = "7/20/17"
$newDate = "4/16/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.TextRange.Text -eq $oldDate) {
                $tf.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout off the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master's date placeholder.
$notesMaster = $p.NotesMaster
for ($ni = 1; $ni -le $notesMaster.Shapes.Count; $ni++) {
    $nshp = $notesMaster.Shapes.Item($ni)
    if ($nshp.HasTextFrame) {
        if ($nshp.TextFrame.TextRange.Text -eq $oldDate) {
            $nshp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Remove the now-obsolete UndoRedoStack diagram shapes from slide 1:
# "Rectangle 62" (UndoRedo / Stack), its "Straight Arrow Connector 57",
# and the "1" TextBox 62 label — these are the last three shapes on the
# slide.
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
$shapes.Item($shapes.Count).Delete()
$shapes.Item($shapes.Count).Delete()
$shapes.Item($shapes.Count).Delete()
